$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1700
$ws.Range("J40").Value = 1900
$ws.Range("L40").Value = 1900
$ws.Range("N40").Value = -2250

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2642.8572
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2642.8572
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2642.8572
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3610.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3865.5278
$ws.Range("I76").Value = 3980.6128
$ws.Range("J76").Value = 3152
$ws.Range("K76").Value = 3980.6128
$ws.Range("L76").Value = 3152
$ws.Range("M76").Value = -3665.6128
$ws.Range("N76").Value = -3782

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3865.5278
$ws.Range("I79").Value = 3980.6128
$ws.Range("J79").Value = 3152
$ws.Range("K79").Value = 3980.6128
$ws.Range("L79").Value = 3152
$ws.Range("M79").Value = -2888.6128
$ws.Range("N79").Value = -5336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1284.5652
$ws.Range("I137").Value = 958.44446
$ws.Range("J137").Value = 2458.6
$ws.Range("K137").Value = 2875.33338
$ws.Range("L137").Value = 7375.799999999999
$ws.Range("M137").Value = -325.33338
$ws.Range("N137").Value = -12475.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3721.44
$ws.Range("I138").Value = 3937.8
$ws.Range("J138").Value = 3577.2
$ws.Range("K138").Value = 11813.4
$ws.Range("L138").Value = 10731.6
$ws.Range("M138").Value = -6673.400000000001
$ws.Range("N138").Value = -21011.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5292860
$ws.Range("I61").Value = 7937558.5
$ws.Range("J61").Value = 3463.476
$ws.Range("K61").Value = 7937558.5
$ws.Range("L61").Value = 3463.476
$ws.Range("M61").Value = -7937346.5
$ws.Range("N61").Value = -3887.476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1078.1072
$ws.Range("I74").Value = 749.3
$ws.Range("J74").Value = 1260.7778
$ws.Range("K74").Value = 749.3
$ws.Range("L74").Value = 1260.7778
$ws.Range("M74").Value = 124.7
$ws.Range("N74").Value = -3008.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1078.1072
$ws.Range("I77").Value = 749.3
$ws.Range("J77").Value = 1260.7778
$ws.Range("K77").Value = 3746.5
$ws.Range("L77").Value = 6303.889
$ws.Range("M77").Value = 621.5
$ws.Range("N77").Value = -15039.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5292860
$ws.Range("I136").Value = 7937558.5
$ws.Range("J136").Value = 3463.476
$ws.Range("K136").Value = 23812675.5
$ws.Range("L136").Value = 10390.428
$ws.Range("M136").Value = -23810125.5
$ws.Range("N136").Value = -15490.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2920.2812
$ws.Range("I134").Value = 2975.2778
$ws.Range("J134").Value = 2849.5715
$ws.Range("K134").Value = 8925.8334
$ws.Range("L134").Value = 8548.7145
$ws.Range("M134").Value = -6390.8334
$ws.Range("N134").Value = -13618.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4514.25
$ws.Range("I31").Value = 987
$ws.Range("J31").Value = 10116.353
$ws.Range("K31").Value = 987
$ws.Range("L31").Value = 10116.353
$ws.Range("M31").Value = -692
$ws.Range("N31").Value = -10706.353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4514.25
$ws.Range("I34").Value = 987
$ws.Range("J34").Value = 10116.353
$ws.Range("K34").Value = 987
$ws.Range("L34").Value = 10116.353
$ws.Range("M34").Value = -785
$ws.Range("N34").Value = -10520.353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2326.9546
$ws.Range("I58").Value = 2615.647
$ws.Range("J58").Value = 1345.4
$ws.Range("K58").Value = 2615.647
$ws.Range("L58").Value = 1345.4
$ws.Range("M58").Value = -2412.647
$ws.Range("N58").Value = -1751.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2326.9546
$ws.Range("I136").Value = 2615.647
$ws.Range("J136").Value = 1345.4
$ws.Range("K136").Value = 7846.941
$ws.Range("L136").Value = 4036.2
$ws.Range("M136").Value = -5296.941
$ws.Range("N136").Value = -9136.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1207.75
$ws.Range("I68").Value = 926.875
$ws.Range("J68").Value = 1432.45
$ws.Range("K68").Value = 2780.625
$ws.Range("L68").Value = 4297.35
$ws.Range("M68").Value = -1969.625
$ws.Range("N68").Value = -5919.35

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1207.75
$ws.Range("I71").Value = 926.875
$ws.Range("J71").Value = 1432.45
$ws.Range("K71").Value = 8341.875
$ws.Range("L71").Value = 12892.05
$ws.Range("M71").Value = -4285.875
$ws.Range("N71").Value = -21004.05

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1906.2609
$ws.Range("I107").Value = 237.1875
$ws.Range("J107").Value = 2796.4333
$ws.Range("K107").Value = 711.5625
$ws.Range("L107").Value = 8389.2999
$ws.Range("M107").Value = 1208.4375
$ws.Range("N107").Value = -12229.2999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 46250
$ws.Range("J69").Value = 46250
$ws.Range("L69").Value = 46250
$ws.Range("N69").Value = -47748

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H72").Value = 46250
$ws.Range("J72").Value = 46250
$ws.Range("L72").Value = 138750
$ws.Range("N72").Value = -146238

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2702.3333
$ws.Range("I122").Value = 1553.5
$ws.Range("K122").Value = 4660.5
$ws.Range("M122").Value = -2210.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2986.3572
$ws.Range("I132").Value = 2001.8334
$ws.Range("J132").Value = 3724.75
$ws.Range("K132").Value = 6005.5002
$ws.Range("L132").Value = 11174.25
$ws.Range("M132").Value = -3475.5002
$ws.Range("N132").Value = -16234.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2558.4285
$ws.Range("I100").Value = 2363
$ws.Range("J100").Value = 2773.4
$ws.Range("K100").Value = 2363
$ws.Range("L100").Value = 2773.4
$ws.Range("M100").Value = -1822
$ws.Range("N100").Value = -3855.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3280.7761
$ws.Range("I132").Value = 3411.125
$ws.Range("J132").Value = 3161.6
$ws.Range("K132").Value = 10233.375
$ws.Range("L132").Value = 9484.799999999999
$ws.Range("M132").Value = -7703.375
$ws.Range("N132").Value = -14544.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1464.1428
$ws.Range("I136").Value = 1462.7059
$ws.Range("K136").Value = 4388.1177
$ws.Range("M136").Value = -1838.1177

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9160
$ws.Range("J41").Value = 9160
$ws.Range("L41").Value = 9160
$ws.Range("N41").Value = -9940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5403399.5
$ws.Range("I132").Value = 2058.8647
$ws.Range("J132").Value = 17159258
$ws.Range("K132").Value = 6176.5941
$ws.Range("L132").Value = 51477774
$ws.Range("M132").Value = -3646.5941
$ws.Range("N132").Value = -51482834
